$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the search-request text from the "Продавец" column (F) to the
# "Поисковый запрос" column (G), updating its content to an actual query.
$ws.Range("F2").ClearContents()
$ws.Range("G2").Value = "ежовик гребенчатый"

# Row 2 no longer needs the extra height required to show the long URL.
$ws.Rows(2).RowHeight = 15.75

# Update the active selection to C2.
$ws.Range("C2").Select()
